# SBDL_ER.xlsx - "Terminando licitacion de requisitos"
# Adds the remaining requirement rows (RQ15-RQ21) to the ITERACION1 sheet,
# resizes the rows that now hold wrapped text, and leaves the workbook
# with the "Prototipos" sheet active/selected.

$wb = $excel.ActiveWorkbook
$iter1 = $wb.Worksheets.Item("ITERACION1")
$prototipos = $wb.Worksheets.Item("Prototipos")

# Make sure we are editing ITERACION1 while we fill in the new rows.
$iter1.Activate()

# --- New requisito rows -----------------------------------------------
# Row 24: RQ15 - Actualizacion / Usuario
$iter1.Range("B24").Value = "RQ15"
$iter1.Range("C24").Value = "Actualizacion"
$iter1.Range("D24").Value = "Usuario"
$iter1.Range("E24").Value = "Actualizar datos del docente"
$iter1.Range("F24").Value = "El sistema debe permitir al usuario docente actualizar sus datos personales, tales como biografia, nombres y apellidos, etc."

# Row 27: RQ16 - Agregar curso al docente / Docente
$iter1.Range("B27").Value = "RQ16"
$iter1.Range("C27").Value = "Agregar"
$iter1.Range("D27").Value = "Docente"
$iter1.Range("E27").Value = "Agregar curso al docente"
$iter1.Range("F27").Value = "El sistema debe permitir al usuario docente agregar el curso que desee dictar a sus alumnos. El curso debe tener todos los detalles que pide el sistema, como categoria, nivel, etc."

# Row 30: RQ17 - Eliminar curso al docente / Docente
$iter1.Range("B30").Value = "RQ17"
$iter1.Range("C30").Value = "Eliminar"
$iter1.Range("D30").Value = "Docente"
$iter1.Range("E30").Value = "Eliminar curso al docente"
$iter1.Range("F30").Value = "El sistema debe permitir al usuario docente eliminar el curso que desee dictar a sus alumnos."

# Row 33: RQ18 - Agrega curso favorito / Alumno
$iter1.Range("B33").Value = "RQ18"
$iter1.Range("C33").Value = "Agregar"
$iter1.Range("D33").Value = "Alumno"
$iter1.Range("E33").Value = "Agrega curso favorito "
$iter1.Range("F33").Value = "El sistema debe permitir al usuario docente agregar un curso a la lista de favoritos quer posee un alumno en su perfil."

# Row 36: RQ19 - Eliminar curso favorito / Alumno
$iter1.Range("B36").Value = "RQ19"
$iter1.Range("C36").Value = "Eliminar"
$iter1.Range("D36").Value = "Alumno"
$iter1.Range("E36").Value = "Eliminar curso favorito"
$iter1.Range("F36").Value = "El sistema debe permitir al usuario alumno eliminar un curso de la lista de favoritos que tenga."

# Row 39: RQ20 - Agregar cursos buscados al historial / Alumno
$iter1.Range("B39").Value = "RQ20"
$iter1.Range("C39").Value = "Agregar"
$iter1.Range("D39").Value = "Alumno"
$iter1.Range("E39").Value = "Agregar cursos buscados al historial"
$iter1.Range("F39").Value = "Se debe generar un historial de busqueda a traves del tiempo del usuario alumno"

# Row 42: RQ21 - Eliminar historial de busqueda / Alumno
$iter1.Range("B42").Value = "RQ21"
$iter1.Range("C42").Value = "Eliminar"
$iter1.Range("D42").Value = "Alumno"
$iter1.Range("E42").Value = "Eliminar historial de busqueda"
$iter1.Range("F42").Value = "El sistema debe permitir al usuario alumno poder eliminar el historial de busqueda que se le muestra al inicio."

# --- Row heights for the newly-filled (wrapped) rows --------------------
$iter1.Rows.Item(25).RowHeight = 12.75
$iter1.Rows.Item(26).RowHeight = 32.25
$iter1.Rows.Item(28).RowHeight = 30
$iter1.Rows.Item(29).RowHeight = 23.25
$iter1.Rows.Item(32).RowHeight = 47.25
$iter1.Rows.Item(33).RowHeight = 42
$iter1.Rows.Item(38).RowHeight = 44.25
$iter1.Rows.Item(41).RowHeight = 36.75
$iter1.Rows.Item(44).RowHeight = 41.25

# --- View state: scroll ITERACION1 so row 29 is at the top and select
# the freshly completed block, then leave "Prototipos" as the active tab.
$iter1.Range("B24:F44").Select()
$excel.ActiveWindow.ScrollRow = 29
$excel.ActiveWindow.ScrollColumn = 1

$prototipos.Activate()
